$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 93 (before the former row 94), shifting
# the old rows 94-105 down to 96-107.
$ws.Range("A94:A95").EntireRow.Insert()

# --- New row 94 ---
$ws.Range("A94").Value = 5
$ws.Range("B94").Value = "Macroferia Regional de Talca"
$ws.Range("C94").Value = "Maule"
$ws.Range("D94").Value = 45127
$ws.Range("E94").Value = 7
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100104
$ws.Range("H94").Value = "Frutos de pepita"
$ws.Range("I94").Value = 100104003
$ws.Range("J94").Value = "Membrillo"
$ws.Range("K94").Value = "Champion"
$ws.Range("L94").Value = "Especial"
$ws.Range("M94").Value = 180
$ws.Range("N94").Value = 12000
$ws.Range("O94").Value = 12000
$ws.Range("P94").Value = 12000
$ws.Range("Q94").Value = "$/bandeja 18 kilos granel"
$ws.Range("R94").Value = "Región de O'Higgins"
$ws.Range("S94").Value = 667
$ws.Range("T94").Value = 18

# --- New row 95 ---
$ws.Range("A95").Value = 5
$ws.Range("B95").Value = "Macroferia Regional de Talca"
$ws.Range("C95").Value = "Maule"
$ws.Range("D95").Value = 45127
$ws.Range("E95").Value = 7
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100104
$ws.Range("H95").Value = "Frutos de pepita"
$ws.Range("I95").Value = 100104003
$ws.Range("J95").Value = "Membrillo"
$ws.Range("K95").Value = "Champion"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 190
$ws.Range("N95").Value = 10000
$ws.Range("O95").Value = 10000
$ws.Range("P95").Value = 10000
$ws.Range("Q95").Value = "$/bandeja 18 kilos granel"
$ws.Range("R95").Value = "Región de O'Higgins"
$ws.Range("S95").Value = 556
$ws.Range("T95").Value = 18
